$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed as a number
# by Excel (losing the original text formatting, e.g. trailing zeros or
# switching to scientific notation) are first forced to Text format.

$ws.Range("D2").Value = '26.574.35'
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("D3").Value = '1.738.97'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.85'
$ws.Range("E5").Value = '  +1.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4922'
$ws.Range("E7").Value = '  +2.61%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2684'
$ws.Range("E8").Value = '  +0.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06289'
$ws.Range("D10").Value = '1.736.32'
$ws.Range("E10").Value = '  +0.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07056'
$ws.Range("E11").Value = '  -0.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.73'
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6147'
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.588'
$ws.Range("E14").Value = '  +0.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '78.15'
$ws.Range("E15").Value = '  +1.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9999'
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '26.591.67'
$ws.Range("E17").Value = '  +0.46%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007313'
$ws.Range("E18").Value = '  +5.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.56'
$ws.Range("E20").Value = '  -1.41%  '
$ws.Range("D21").Value = '1.959.99'
$ws.Range("E21").Value = '  +0.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.597'
$ws.Range("E22").Value = '  +0.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.726'
$ws.Range("E23").Value = '  -2.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.258'
$ws.Range("E24").Value = '  -1.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.33'
$ws.Range("E25").Value = '  +2.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.47'
$ws.Range("E26").Value = '  +0.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.419'
$ws.Range("E27").Value = '  +0.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.763'
$ws.Range("E28").Value = '  -1.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '107.68'
$ws.Range("E29").Value = '  +0.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.049'
$ws.Range("E30").Value = '  +1.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08053'
$ws.Range("E31").Value = '  +0.67%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.733'
$ws.Range("E32").Value = '  -0.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04622'
$ws.Range("E33").Value = '  +1.31%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.612'
$ws.Range("E34").Value = '  -0.09%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.019'
$ws.Range("E35").Value = '  +2.83%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6391'
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9000'
$ws.Range("E37").Value = '  -4.43%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.046'
$ws.Range("E38").Value = '  +2.62%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.404'
$ws.Range("E39").Value = '  -0.20%  '
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.003'
$ws.Range("E40").Value = '  -0.34%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01503'
$ws.Range("E41").Value = '  +0.19%  '
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.97'
$ws.Range("E42").Value = '  -5.45%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.418'
$ws.Range("E43").Value = '  -4.13%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3927'
$ws.Range("E44").Value = '  +0.52%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.910'
$ws.Range("E45").Value = '  -0.58%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1186'
$ws.Range("E46").Value = '  -0.56%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05399'
$ws.Range("E47").Value = '  +1.39%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '30.59'
$ws.Range("E48").Value = '  -0.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.791'
$ws.Range("E49").Value = '  -1.32%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.271'
$ws.Range("E50").Value = '  -0.08%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.73'
$ws.Range("E51").Value = '  +0.72%  '
